$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.849.85'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.13%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.468.81'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.62%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '568.03'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.85%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '166.57'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.33%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.513'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.69%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.176'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +12.59%  '
$ws.Range('E10').Value = '  -1.49%  '
$ws.Range('E11').Value = '  +2.28%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.67'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.81%  '
$ws.Range('E13').Value = '  +8.34%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '69.713.55'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.921.43'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.54%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '23.90'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.45%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.461.47'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.30%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.84'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.48%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '343.57'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.08%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.17'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +4.44%  '
$ws.Range('E21').Value = '  +3.29%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.02'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +8.37%  '
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '66.26'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.48%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.90'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +6.63%  '
$ws.Range('B26').NumberFormat = '@'
$ws.Range('B26').Value = 'WrappedeETH'
$ws.Range('B26').Style = 'Normal'
$ws.Range('C26').NumberFormat = '@'
$ws.Range('C26').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('C26').Style = 'Normal'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.584.36'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.25%  '
$ws.Range('B27').NumberFormat = '@'
$ws.Range('B27').Value = 'Aptos'
$ws.Range('B27').Style = 'Normal'
$ws.Range('C27').NumberFormat = '@'
$ws.Range('C27').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('C27').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.55'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +5.56%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.995'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.70%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0853'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.84%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.33'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.18%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.26'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +10.92%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '450.07'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +6.63%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.05%  '
$ws.Range('E34').Value = '  +1.16%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '160.28'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.53%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '19.05'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.41%  '
$ws.Range('B37').NumberFormat = '@'
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('B37').Style = 'Normal'
$ws.Range('C37').NumberFormat = '@'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('C37').Style = 'Normal'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.110'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.37%  '
$ws.Range('B38').NumberFormat = '@'
$ws.Range('B38').Value = 'USDe'
$ws.Range('B38').Style = 'Normal'
$ws.Range('C38').NumberFormat = '@'
$ws.Range('C38').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('C38').Style = 'Normal'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.01%  '
$ws.Range('E39').Value = '  +2.49%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.307'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.03%  '
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('B41').Style = 'Normal'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('C41').Style = 'Normal'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.49'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.87%  '
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'Stacks'
$ws.Range('B42').Style = 'Normal'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('C42').Style = 'Normal'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.55'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +5.49%  '
$ws.Range('E43').Value = '  +3.61%  '
$ws.Range('E44').Value = '  +6.93%  '
$ws.Range('E45').Value = '  +1.39%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '132.83'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.61%  '
$ws.Range('E47').Value = '  +0.76%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.492'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.47%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.564'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.48%  '
$ws.Range('E50').Value = '  +1.41%  '
$ws.Range('E51').Value = '  +2.72%  '
